$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.472621202468872
$ws.Range("B1").Value = 2.299251317977905
$ws.Range("C1").Value = 5.161932468414307
$ws.Range("D1").Value = 3.313603639602661
$ws.Range("E1").Value = 1.085882782936096
